$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing row 51 values before they get overwritten, so we
# can re-insert them (unchanged) as the new row 52.
$oldVals = @{}
for ($col = 1; $col -le 18; $col++) {
    $oldVals[$col] = $ws.Cells.Item(51, $col).Value()
}
$oldNumFmt = $ws.Cells.Item(51, 4).NumberFormat

# Insert a new blank row at 52 - this pushes the former row 52
# (date 2021-07-23 / 44400 ...) down to row 53, matching the diff.
$ws.Rows.Item(52).EntireRow.Insert()

# Re-populate new row 52 with the values row 51 used to hold.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(52, $col).Value = $oldVals[$col]
}
$ws.Cells.Item(52, 4).NumberFormat = $oldNumFmt

# Now update row 51 in place with the new observation's data.
$ws.Cells.Item(51, 4).Value = 44448
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 12000
$ws.Cells.Item(51, 12).Value = 13000
$ws.Cells.Item(51, 13).Value = 12450
$ws.Cells.Item(51, 16).Value = 498
